# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-13 20:18:55
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet lists the
# recorder names/emails for a session. For a specific set of rows the order of the
# two comma-separated entries was flipped from "System, dnasr281@gmail.com" to
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

# Exact rows (column G) whose "Recorded By" text had its two parts swapped.
$targetRows = @(
    2,3,4,5,6,7,
    16,17,
    22,23,
    37,38,
    43,44,
    58,59,
    64,65,
    79,80,
    85,86,87,88,89,90,
    99,100,
    105,106,107,108,109,110,
    119,120,
    125,126,127,128,129,130,
    139,140,
    145,146,147,148,149,150,
    159,160,
    165,166,167,168,169,170,
    179,180,
    185,186,
    200,201,
    206,207,
    221,222,
    227,228,
    242,243
)

foreach ($r in $targetRows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
